$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while forcing text storage so that
# numeric-looking strings (e.g. "1.000", "0.9992") are preserved verbatim
# instead of being normalized into numbers by Excel.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.548.48"
Set-TextValue $ws.Range("E2") "  -0.15%  "
Set-TextValue $ws.Range("D3") "1.922.29"
Set-TextValue $ws.Range("E3") "  -0.09%  "
Set-TextValue $ws.Range("D4") "1.000"
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "245.35"
Set-TextValue $ws.Range("E5") "  -0.88%  "
Set-TextValue $ws.Range("D6") "0.9992"
Set-TextValue $ws.Range("D7") "0.4870"
Set-TextValue $ws.Range("E7") "  +2.72%  "
Set-TextValue $ws.Range("D8") "0.2905"
Set-TextValue $ws.Range("E8") "  -0.81%  "
Set-TextValue $ws.Range("D9") "0.06723"
Set-TextValue $ws.Range("E9") "  -1.28%  "
Set-TextValue $ws.Range("D10") "110.79"
Set-TextValue $ws.Range("E10") "  +4.82%  "
Set-TextValue $ws.Range("D11") "19.11"
Set-TextValue $ws.Range("E11") "  +3.40%  "
Set-TextValue $ws.Range("D12") "1.914.72"
Set-TextValue $ws.Range("E12") "  -0.42%  "
Set-TextValue $ws.Range("D13") "0.07576"
Set-TextValue $ws.Range("E13") "  -1.97%  "
Set-TextValue $ws.Range("D14") "5.290"
Set-TextValue $ws.Range("E14") "  -0.56%  "
Set-TextValue $ws.Range("D15") "0.6699"
Set-TextValue $ws.Range("E15") "  -0.41%  "
Set-TextValue $ws.Range("D16") "296.14"
Set-TextValue $ws.Range("E16") "  +2.95%  "
Set-TextValue $ws.Range("D17") "30.541.99"
Set-TextValue $ws.Range("E17") "  -0.28%  "
Set-TextValue $ws.Range("E18") "  +0.46%  "
Set-TextValue $ws.Range("D19") "0.9990"
Set-TextValue $ws.Range("E19") "  -0.14%  "
Set-TextValue $ws.Range("D20") "0.000007576"
Set-TextValue $ws.Range("E20") "  -0.88%  "
Set-TextValue $ws.Range("D21") "5.557"
Set-TextValue $ws.Range("E21") "  +2.10%  "
Set-TextValue $ws.Range("D22") "2.163.87"
Set-TextValue $ws.Range("E22") "  +0.09%  "
Set-TextValue $ws.Range("D23") "1.000"
Set-TextValue $ws.Range("E23") "  +0.04%  "
Set-TextValue $ws.Range("D24") "6.456"
Set-TextValue $ws.Range("E24") "  +2.13%  "
Set-TextValue $ws.Range("D25") "9.464"
Set-TextValue $ws.Range("E25") "  +0.43%  "
Set-TextValue $ws.Range("D26") "164.70"
Set-TextValue $ws.Range("E26") "  -2.14%  "
Set-TextValue $ws.Range("D27") "20.26"
Set-TextValue $ws.Range("E27") "  -3.09%  "
Set-TextValue $ws.Range("D28") "2.105"
Set-TextValue $ws.Range("E28") "  -1.74%  "
Set-TextValue $ws.Range("D29") "0.1073"
Set-TextValue $ws.Range("E29") "  -1.00%  "
Set-TextValue $ws.Range("D30") "1.444"
Set-TextValue $ws.Range("E30") "  +5.65%  "
Set-TextValue $ws.Range("D31") "4.155"
Set-TextValue $ws.Range("E31") "  -1.22%  "
Set-TextValue $ws.Range("D32") "4.060"
Set-TextValue $ws.Range("E32") "  -1.78%  "
Set-TextValue $ws.Range("D33") "0.05030"
Set-TextValue $ws.Range("E33") "  -0.47%  "
Set-TextValue $ws.Range("D34") "0.7407"
Set-TextValue $ws.Range("E34") "  -0.45%  "
Set-TextValue $ws.Range("D35") "1.139"
Set-TextValue $ws.Range("E35") "  -1.86%  "
Set-TextValue $ws.Range("D36") "0.9993"
Set-TextValue $ws.Range("E36") "  -0.01%  "
Set-TextValue $ws.Range("D37") "2.709"
Set-TextValue $ws.Range("E37") "  -1.27%  "
Set-TextValue $ws.Range("D38") "0.02027"
Set-TextValue $ws.Range("E38") "  -2.60%  "
Set-TextValue $ws.Range("D39") "2.684"
Set-TextValue $ws.Range("E39") "  -0.24%  "
Set-TextValue $ws.Range("D40") "110.44"
Set-TextValue $ws.Range("E40") "  -0.93%  "
Set-TextValue $ws.Range("D41") "2.021"
Set-TextValue $ws.Range("E41") "  -2.50%  "
Set-TextValue $ws.Range("D42") "0.4436"
Set-TextValue $ws.Range("E42") "  +0.84%  "
Set-TextValue $ws.Range("D43") "0.8667"
Set-TextValue $ws.Range("E43") "  -1.77%  "
Set-TextValue $ws.Range("D44") "70.98"
Set-TextValue $ws.Range("E44") "  +5.05%  "
Set-TextValue $ws.Range("D45") "5.831"
Set-TextValue $ws.Range("E45") "  -2.25%  "
Set-TextValue $ws.Range("D46") "0.9987"
Set-TextValue $ws.Range("E46") "  -0.13%  "
Set-TextValue $ws.Range("D47") "7.265"
Set-TextValue $ws.Range("E47") "  -0.24%  "
Set-TextValue $ws.Range("D48") "48.62"
Set-TextValue $ws.Range("E48") "  +0.73%  "
Set-TextValue $ws.Range("D51") "0.2534"
Set-TextValue $ws.Range("E51") "  +3.11%  "

# Rows 49 and 50 swap places (Algorand <-> EnergySwap) with updated data.
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "9.257"
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D50") "0.1233"
$ws.Range("E50").Value = "  -0.41%  "
